$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (copy formatting from the previous header cell, then set text)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Mex_US_Border"

# New column data (rows 2-41)
$hValues = @(
    5.5,
    6.5,
    5.4,
    5.75,
    6.25,
    33.4,
    8,
    10.4,
    12.75,
    10.75,
    8.8,
    8,
    9.5,
    7.4,
    6.75,
    9,
    9.5,
    14.5,
    10,
    8.6,
    11.5,
    10.5,
    8.4,
    9.5,
    9,
    7.25,
    6,
    7.8,
    8.75,
    10.5,
    12,
    11.2,
    9,
    8.2,
    9.75,
    9,
    9,
    8.5,
    10.25,
    11
)

for ($i = 0; $i -lt $hValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $hValues[$i]
}

# Match the author's final selection state (entire column H selected, active cell H1)
$ws.Range("H1:H1048576").Select()
